# edit.ps1
# Applies the data changes described by the commit to DOC/flash分配.xlsx
# (sheet "Sheet1"):
#   - clear the "出厂值" (factory/default value) column G for rows 2-9
#     (the fast/slow cfar threshold rows), since those constants are no
#     longer fixed defaults,
#   - bump the G13 "UPSSA0" base-address related value from 3800 to 4000,
#   - bump the G14 delay/threshold value from 32 to 10 (T=4.8s slow check
#     fix),
#   - leave the cursor / active selection on G7, matching where the user
#     was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously hard-coded factory values in G2:G9 (keep cell
# formatting/style, only remove the stored number).
$ws.Range("G2:G9").ClearContents()

# Update the two tuned numeric constants.
$ws.Range("G13").Value = 4000
$ws.Range("G14").Value = 10

# Restore/move the active selection to G7.
$ws.Range("G7").Select()
